$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) cells that change: force text storage so Excel
# doesn't auto-coerce numeric-looking strings into floating point numbers
# (they are plain text cells in the source workbook, t="inlineStr").
$dCells = @{
    "D2"  = "33.633.39"
    "D3"  = "1.769.66"
    "D5"  = "223.61"
    "D8"  = "31.72"
    "D9"  = "0.288"
    "D10" = "0.0684"
    "D11" = "0.0935"
    "D13" = "10.99"
    "D14" = "1.755.12"
    "D15" = "33.670.91"
    "D18" = "66.34"
    "D20" = "237.58"
    "D22" = "10.51"
    "D23" = "3.99"
    "D25" = "159.32"
    "D26" = "16.04"
    "D31" = "0.0509"
    "D35" = "1.379.09"
    "D38" = "0.0184"
    "D41" = "77.59"
    "D42" = "2.66"
    "D43" = "0.901"
    "D44" = "13.46"
    "D48" = "106.83"
    "D50" = "1.923.99"
    "D51" = "1.00"
}

foreach ($addr in $dCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $dCells.Keys) {
    $ws.Range($addr).Value = $dCells[$addr]
}
foreach ($addr in $dCells.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# --- Column B / C (swap ARBITRUM <-> MXToken between rows 42 and 43)
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

# --- Column E (Volume 1h) cells that change: plain text percentage values
$ws.Range("E2").Value  = "  -0.92%  "
$ws.Range("E3").Value  = "  -0.91%  "
$ws.Range("E5").Value  = "  +0.84%  "
$ws.Range("E6").Value  = "  -1.02%  "
$ws.Range("E7").Value  = "  +0.03%  "
$ws.Range("E8").Value  = "  +0.95%  "
$ws.Range("E9").Value  = "  +1.46%  "
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  +4.12%  "
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("E44").Value = "  +13.86%  "
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("E46").Value = "  +14.54%  "
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  +0.24%  "
